$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column A width (wide column for long ngram list text)
$ws.Columns.Item(1).ColumnWidth = 103

# Populate column H (ngram(1,3) after preprocessing) for data rows 2-43 first
$ws.Range("H2").Value = '[''0.7'', ''0.59'', ''0.58'']'
$ws.Range("H3").Value = '[''0.77'', ''0.73'', ''0.71'']'
$ws.Range("H4").Value = '[''0.38'', ''0.38'', ''0.38'']'
$ws.Range("H5").Value = '[''0.49'', ''0.43'', ''0.38'']'
$ws.Range("H6").Value = '[''0.48'', ''0.46'', ''0.43'']'
$ws.Range("H7").Value = '[''0.27'', ''0.26'', ''0.25'']'
$ws.Range("H8").Value = '[''0.26'', ''0.25'', ''0.24'']'
$ws.Range("H9").Value = '[''0.67'', ''0.37'', ''0.32'']'
$ws.Range("H10").Value = '[''0.48'', ''0.44'', ''0.26'']'
$ws.Range("H11").Value = '[''0.43'', ''0.43'', ''0.38'']'
$ws.Range("H12").Value = '[''0.75'', ''0.47'', ''0.43'']'
$ws.Range("H13").Value = '[''0.77'', ''0.73'', ''0.71'']'
$ws.Range("H14").Value = '[''0.43'', ''0.41'', ''0.38'']'
$ws.Range("H15").Value = '[''0.71'', ''0.4'', ''0.36'']'
$ws.Range("H16").Value = '[''0.62'', ''0.51'', ''0.5'']'
$ws.Range("H17").Value = '[''0.49'', ''0.46'', ''0.46'']'
$ws.Range("H18").Value = '[''0.74'', ''0.52'', ''0.51'']'
$ws.Range("H19").Value = '[''0.33'', ''0.29'', ''0.26'']'
$ws.Range("H20").Value = '[''0.5'', ''0.43'', ''0.4'']'
$ws.Range("H21").Value = '[''0.48'', ''0.45'', ''0.41'']'
$ws.Range("H22").Value = '[''0.94'', ''0.68'', ''0.65'']'
$ws.Range("H23").Value = '[''0.46'', ''0.46'', ''0.45'']'
$ws.Range("H24").Value = '[''0.5'', ''0.37'', ''0.33'']'
$ws.Range("H25").Value = '[''0.76'', ''0.57'', ''0.55'']'
$ws.Range("H26").Value = '[''0.58'', ''0.58'', ''0.56'']'
$ws.Range("H27").Value = '[''0.48'', ''0.47'', ''0.39'']'
$ws.Range("H28").Value = '[''0.84'', ''0.67'', ''0.62'']'
$ws.Range("H29").Value = '[''0.58'', ''0.57'', ''0.5'']'
$ws.Range("H30").Value = '[''0.81'', ''0.74'', ''0.65'']'
$ws.Range("H31").Value = '[''0.43'', ''0.34'', ''0.33'']'
$ws.Range("H32").Value = '[''0.49'', ''0.46'', ''0.45'']'
$ws.Range("H33").Value = '[''0.72'', ''0.46'', ''0.45'']'
$ws.Range("H34").Value = '[''0.7'', ''0.53'', ''0.48'']'
$ws.Range("H35").Value = '[''0.44'', ''0.4'', ''0.34'']'
$ws.Range("H36").Value = '[''0.84'', ''0.64'', ''0.6'']'
$ws.Range("H37").Value = '[''0.5'', ''0.33'', ''0.3'']'
$ws.Range("H38").Value = '[''0.36'', ''0.35'', ''0.33'']'
$ws.Range("H39").Value = '[''0.34'', ''0.28'', ''0.25'']'
$ws.Range("H40").Value = '[''0.4'', ''0.36'', ''0.32'']'
$ws.Range("H41").Value = '[''0.36'', ''0.36'', ''0.34'']'
$ws.Range("H42").Value = '[''0.43'', ''0.39'', ''0.39'']'
$ws.Range("H43").Value = '[''0.46'', ''0.45'', ''0.38'']'

# Populate column K (ngram(1,3) without preprocessing) for data rows 2-43 next
$ws.Range("K2").Value = '[''0.55'', ''0.45'', ''0.42'']'
$ws.Range("K3").Value = '[''0.67'', ''0.61'', ''0.56'']'
$ws.Range("K4").Value = '[''0.4'', ''0.39'', ''0.37'']'
$ws.Range("K5").Value = '[''0.51'', ''0.47'', ''0.46'']'
$ws.Range("K6").Value = '[''0.53'', ''0.47'', ''0.44'']'
$ws.Range("K7").Value = '[''0.36'', ''0.32'', ''0.31'']'
$ws.Range("K8").Value = '[''0.22'', ''0.22'', ''0.21'']'
$ws.Range("K9").Value = '[''0.66'', ''0.36'', ''0.35'']'
$ws.Range("K10").Value = '[''0.5'', ''0.43'', ''0.24'']'
$ws.Range("K11").Value = '[''0.48'', ''0.44'', ''0.39'']'
$ws.Range("K12").Value = '[''0.69'', ''0.4'', ''0.38'']'
$ws.Range("K13").Value = '[''0.62'', ''0.55'', ''0.54'']'
$ws.Range("K14").Value = '[''0.59'', ''0.56'', ''0.53'']'
$ws.Range("K15").Value = '[''0.62'', ''0.33'', ''0.28'']'
$ws.Range("K16").Value = '[''0.54'', ''0.53'', ''0.5'']'
$ws.Range("K17").Value = '[''0.53'', ''0.47'', ''0.44'']'
$ws.Range("K18").Value = '[''0.68'', ''0.46'', ''0.45'']'
$ws.Range("K19").Value = '[''0.17'', ''0.16'', ''0.16'']'
$ws.Range("K20").Value = '[''0.52'', ''0.41'', ''0.37'']'
$ws.Range("K21").Value = '[''0.45'', ''0.45'', ''0.45'']'
$ws.Range("K22").Value = '[''0.82'', ''0.58'', ''0.57'']'
$ws.Range("K23").Value = '[''0.45'', ''0.44'', ''0.43'']'
$ws.Range("K24").Value = '[''0.61'', ''0.41'', ''0.4'']'
$ws.Range("K25").Value = '[''0.73'', ''0.57'', ''0.5'']'
$ws.Range("K26").Value = '[''0.44'', ''0.42'', ''0.41'']'
$ws.Range("K27").Value = '[''0.55'', ''0.46'', ''0.46'']'
$ws.Range("K28").Value = '[''0.83'', ''0.66'', ''0.65'']'
$ws.Range("K29").Value = '[''0.62'', ''0.58'', ''0.55'']'
$ws.Range("K30").Value = '[''0.7'', ''0.62'', ''0.57'']'
$ws.Range("K31").Value = '[''0.4'', ''0.39'', ''0.34'']'
$ws.Range("K32").Value = '[''0.47'', ''0.45'', ''0.45'']'
$ws.Range("K33").Value = '[''0.67'', ''0.45'', ''0.41'']'
$ws.Range("K34").Value = '[''0.63'', ''0.46'', ''0.37'']'
$ws.Range("K35").Value = '[''0.41'', ''0.37'', ''0.33'']'
$ws.Range("K36").Value = '[''0.7'', ''0.48'', ''0.47'']'
$ws.Range("K37").Value = '[''0.5'', ''0.35'', ''0.34'']'
$ws.Range("K38").Value = '[''0.29'', ''0.25'', ''0.24'']'
$ws.Range("K39").Value = '[''0.4'', ''0.33'', ''0.33'']'
$ws.Range("K40").Value = '[''0.41'', ''0.39'', ''0.34'']'
$ws.Range("K41").Value = '[''0.4'', ''0.38'', ''0.32'']'
$ws.Range("K42").Value = '[''0.44'', ''0.42'', ''0.42'']'
$ws.Range("K43").Value = '[''0.4'', ''0.34'', ''0.33'']'

# Finally add the header labels in row 1
$ws.Range("H1").Value = 'ngram(1,3) after prepro'
$ws.Range("K1").Value = 'ngram(1,3) without prepro'

# Scroll the view down and select the last populated cell, matching the saved view state
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("A43").Select() | Out-Null
